# Apply the "Setting prefix explicitly in excelparser" edit to onto.xlsx
#
# Summary of the change:
#  - ImportedOntologies sheet gains a 4th documented column ("base_iri_root"
#    moves to C, a new "Comments" column D is introduced) and two example
#    rows of data showing the new explicit-prefix feature.
#  - Concepts sheet: the SpecialMolecule / AnotherSpecialMolecule examples
#    are updated to use the "emmo:" prefix (instead of the old
#    "emmo-inferred-chemistry2:" prefix) and two new example rows
#    (ANewTestClass / AnotherNewTestClass) are appended.
#  - The workbook is left with the Concepts sheet as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. ImportedOntologies sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ImportedOntologies")

# New "Comments" example cell in the header/description block (row 2),
# styled like the existing grey/italic description cells but without a
# border (re-using font 2 / fill 4 from the existing style palette).
$d2 = $ws2.Range("D2")
$d2.Value2 = "Comment for test"
$d2.Font.Italic = $true
$d2.Interior.Color = 13553360
$d2.Interior.PatternColor = 15189684
$d2.WrapText = $true

# Row 3: existing "emmo" import line gets example base_iri_root / comment
$ws2.Range("C3").Value2 = "  "
$ws2.Range("D3").Value2 = "Do not give base_iri_root, but leave spaces there (so it is not completely empty)"

# Row 4: new example import with explicit local prefix + base_iri_root
$ws2.Range("A4").Value2 = "imported_onto/ontology.ttl"
$ws2.Range("B4").Value2 = "testonto"
$ws2.Range("C4").Value2 = "http://ontology.info/"

# ---------------------------------------------------------------------
# 2. Concepts sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Concepts")

# Row 25 (SpecialMolecule): subClassOf / relation now use the "emmo:"
# prefix, and the comment explains the prefix is chosen via
# ImportedOntologies rather than being a known limitation.
$ws3.Range("G25").Value2 = "emmo:Molecule"
$ws3.Range("I25").Value2 = "hasPart some Atom"
$ws3.Range("J25").Value2 = "Test giving prefix from emmo for subclass of and relations. Prefix is set by choice in ImportedOntologies."

# Row 26 (AnotherSpecialMolecule): relation now uses "emmo:" prefix too.
$ws3.Range("I26").Value2 = "emmo:hasPart  some emmo:Atom"
$ws3.Range("J26").Value2 = "Test giving prefix to relations."

# New row 27: example class picking up the prefix of the first import.
$ws3.Range("A27").Value2 = "ANewTestClass"
$ws3.Range("G27").Value2 = "testonto:TestClass"
$ws3.Range("J27").Value2 = "Check that prefix is set to first import (ontology.ttl)"

# New row 28: example class picking up the prefix of the sub-import.
$ws3.Range("A28").Value2 = "AnotherNewTestClass"
$ws3.Range("G28").Value2 = "testonto:TestClass2"
$ws3.Range("J28").Value2 = "Check that prefix is set to subimport (ontology.ttl)"

# ---------------------------------------------------------------------
# 3. View state: ImportedOntologies selection moves to D4 (no longer the
#    active tab); Concepts becomes the active tab, scrolled/zoomed to
#    show the newly added rows, with H34 selected.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("D4").Select()

$ws3.Activate()
$excel.ActiveWindow.Zoom = 85
$ws3.Range("H34").Select()
